$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at M, pushing the existing "Email" column from M to N.
$ws.Columns.Item(13).Insert()

# New header for the inserted column: the 2020-12-14 "chamada" (attendance) date.
# Format the cell as text first so the yyyy-mm-dd-looking string is not
# auto-converted into a date serial value (it keeps the bold/border style that
# it inherited from the column insert).
$ws.Cells.Item(1, 13).NumberFormat = "@"
$ws.Cells.Item(1, 13).Value = "2020-12-14"

# New per-student attendance values for the 2020-12-14 column.
$ws.Cells.Item(2, 13).Value = 21.62
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(5, 13).Value = 100
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(15, 13).NumberFormat = "@"
$ws.Cells.Item(15, 13).Value = ""
$ws.Cells.Item(16, 13).Value = 0
$ws.Cells.Item(17, 13).NumberFormat = "@"
$ws.Cells.Item(17, 13).Value = ""
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(19, 13).Value = 59.46
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(21, 13).Value = 0
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(26, 13).Value = 0
$ws.Cells.Item(27, 13).Value = 0
$ws.Cells.Item(28, 13).NumberFormat = "@"
$ws.Cells.Item(28, 13).Value = ""
$ws.Cells.Item(29, 13).NumberFormat = "@"
$ws.Cells.Item(29, 13).Value = ""
$ws.Cells.Item(30, 13).Value = 0
$ws.Cells.Item(31, 13).Value = 0
$ws.Cells.Item(32, 13).Value = 0
$ws.Cells.Item(33, 13).Value = 0
$ws.Cells.Item(34, 13).Value = 0
$ws.Cells.Item(35, 13).Value = 0
$ws.Cells.Item(36, 13).Value = 83.78
$ws.Cells.Item(37, 13).Value = 13.51

# Row 31 also had its "Resenha Regime de Metas" (column L) score corrected to 0.
$ws.Cells.Item(31, 12).Value = 0
